# Add the new "acc_std" sheet (accuracy standard-deviation table) in between
# "acc_median" and "MCC", matching the layout of the other accuracy sheets.
$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("acc_median")
$mcc = $wb.Worksheets.Item("MCC")
$template.Copy($mcc)

$ws = $wb.Worksheets.Item("acc_median (2)")
$ws.Name = "acc_std"

# Populate the standard-deviation values for each algorithm / labeling method.
$ws.Range("C3").Value = 0.029214594095834
$ws.Range("D3").Value = 0.035220728260615897
$ws.Range("E3").Value = 0.035063039566093197
$ws.Range("F3").Value = 0.037106497956612701
$ws.Range("G3").Value = 0.030257378733034499
$ws.Range("H3").Value = 0.0322453172684634
$ws.Range("I3").Value = 0.035948071970498603
$ws.Range("J3").Value = 0.0365008378763194

$ws.Range("C4").Value = 0.038291428432796003
$ws.Range("D4").Value = 0.031728323276584001
$ws.Range("E4").Value = 0.037438450379306203
$ws.Range("F4").Value = 0.030124152013660001
$ws.Range("G4").Value = 0.043739668428701903
$ws.Range("H4").Value = 0.055400128667852003
$ws.Range("I4").Value = 0.0346453372450055
$ws.Range("J4").Value = 0.042237548043871899

$ws.Range("C5").Value = 0.040196455347759397
$ws.Range("D5").Value = 0.027863197239903801
$ws.Range("E5").Value = 0.035736133907901299
$ws.Range("F5").Value = 0.039507948534766801
$ws.Range("G5").Value = 0.032470415844113802
$ws.Range("H5").Value = 0.033472857077707398
$ws.Range("I5").Value = 0.032491884601928402
$ws.Range("J5").Value = 0.046857103548961797

$ws.Range("C6").Value = 0.039146928646325102
$ws.Range("D6").Value = 0.037372939184446903
$ws.Range("E6").Value = 0.032603106522415902
$ws.Range("F6").Value = 0.031157684934527901
$ws.Range("G6").Value = 0.033544812548609498
$ws.Range("H6").Value = 0.033075262216712302
$ws.Range("I6").Value = 0.038571732942935102
$ws.Range("J6").Value = 0.033484610142145502

$ws.Range("C7").Value = 0.037845427263483097
$ws.Range("D7").Value = 0.040108531306391401
$ws.Range("E7").Value = 0.034260266767959097
$ws.Range("F7").Value = 0.106050414498165
$ws.Range("G7").Value = 0.0354827247482652
$ws.Range("H7").Value = 0.039303360608536397
$ws.Range("I7").Value = 0.034378994162921103
$ws.Range("J7").Value = 0.037806996728791599

$ws.Range("C8").Value = 0.038728794916282602
$ws.Range("D8").Value = 0.039752416138198902
$ws.Range("E8").Value = 0.037039750476537597
$ws.Range("F8").Value = 0.0567199863580563
$ws.Range("G8").Value = 0.043484037575017301
$ws.Range("H8").Value = 0.036414796197858403
$ws.Range("I8").Value = 0.0386802230799798
$ws.Range("J8").Value = 0.042875929902538697

$ws.Select()
$ws.Range("A2").Select()
